$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.249.92"
$ws.Range("E2").Value = "  +4.30%  "
$ws.Range("D3").Value = "2.964.35"
$ws.Range("E3").Value = "  +2.64%  "
$ws.Range("E4").Value = "  +0.56%  "
$ws.Range("D5").Value = "579.47"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").Value = "152.60"
$ws.Range("E6").Value = "  +6.31%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "2.962.92"
$ws.Range("E8").Value = "  +2.76%  "
$ws.Range("D9").Value = "0.510"
$ws.Range("E9").Value = "  +1.80%  "
$ws.Range("E10").Value = "  +4.11%  "
$ws.Range("D11").Value = "0.152"
$ws.Range("E11").Value = "  +2.14%  "
$ws.Range("D12").Value = "0.445"
$ws.Range("E12").Value = "  +3.23%  "
$ws.Range("D13").Value = "0.0000241"
$ws.Range("E13").Value = "  +3.47%  "
$ws.Range("D14").Value = "34.37"
$ws.Range("E14").Value = "  +7.14%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "3.448.75"
$ws.Range("E16").Value = "  +2.82%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "64.190.02"
$ws.Range("E17").Value = "  +4.65%  "
$ws.Range("D18").Value = "6.88"
$ws.Range("E18").Value = "  +4.39%  "
$ws.Range("D19").Value = "2.959.29"
$ws.Range("E19").Value = "  +2.36%  "
$ws.Range("D20").Value = "445.39"
$ws.Range("E20").Value = "  +2.57%  "
$ws.Range("D21").Value = "13.51"
$ws.Range("E21").Value = "  +2.30%  "
$ws.Range("D22").Value = "0.675"
$ws.Range("E22").Value = "  +3.20%  "
$ws.Range("D23").Value = "7.21"
$ws.Range("E23").Value = "  +4.62%  "
$ws.Range("D24").Value = "80.32"
$ws.Range("E24").Value = "  +1.41%  "
$ws.Range("D25").Value = "10.95"
$ws.Range("E25").Value = "  +8.59%  "
$ws.Range("D26").Value = "12.30"
$ws.Range("E26").Value = "  +3.84%  "
$ws.Range("D27").Value = "2.18"
$ws.Range("E27").Value = "  +7.38%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  +2.32%  "
$ws.Range("D30").Value = "7.63"
$ws.Range("E30").Value = "  +8.98%  "
$ws.Range("D31").Value = "2.17"
$ws.Range("E31").Value = "  +4.78%  "
$ws.Range("D32").Value = "2.55"
$ws.Range("E32").Value = "  +1.82%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "26.54"
$ws.Range("E33").Value = "  +3.68%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.109"
$ws.Range("E34").Value = "  +3.05%  "
$ws.Range("E35").Value = "  +0.48%  "
$ws.Range("D36").Value = "0.973"
$ws.Range("E36").Value = "  +1.53%  "
$ws.Range("D37").Value = "5.63"
$ws.Range("E37").Value = "  +3.64%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "3.08"
$ws.Range("E38").Value = "  +3.70%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "2.10"
$ws.Range("E39").Value = "  +7.76%  "
$ws.Range("D40").Value = "49.06"
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("D41").Value = "43.56"
$ws.Range("E41").Value = "  +13.93%  "
$ws.Range("D42").Value = "0.119"
$ws.Range("E42").Value = "  +2.51%  "
$ws.Range("D43").Value = "0.293"
$ws.Range("E43").Value = "  +9.12%  "
$ws.Range("D44").Value = "8.30"
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("D45").Value = "384.30"
$ws.Range("E45").Value = "  +14.14%  "
$ws.Range("D46").Value = "2.753.28"
$ws.Range("E46").Value = "  +2.73%  "
$ws.Range("D47").Value = "0.0348"
$ws.Range("E47").Value = "  +5.07%  "
$ws.Range("D48").Value = "134.25"
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "0.000223"
$ws.Range("E50").Value = "  +13.02%  "
$ws.Range("E51").Value = "  +2.84%  "
